# Hortaliza, Vega Monumental Concepción - Zanahoria
# Weekly update: insert the newest week's prices at the top of the data
# block (rows 113-114), pushing all existing weekly rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 113 (shifts existing rows 113:128 down to 115:130,
# carrying their formatting - e.g. the date style on column D - along with them).
$ws.Range("A113:A114").EntireRow.Insert()

# New row 113 - Primera quality, newest week
$ws.Cells.Item(113, 1).Value = 11
$ws.Cells.Item(113, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(113, 3).Value = "Bíobío"
$ws.Cells.Item(113, 4).Value = 44461
$ws.Cells.Item(113, 5).Value = 8
$ws.Cells.Item(113, 6).Value = 100114013
$ws.Cells.Item(113, 7).Value = "Zanahoria"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 600
$ws.Cells.Item(113, 11).Value = 8000
$ws.Cells.Item(113, 12).Value = 9000
$ws.Cells.Item(113, 13).Value = 8500
$ws.Cells.Item(113, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(113, 15).Value = "Región de Ñuble"
$ws.Cells.Item(113, 16).Value = 425
$ws.Cells.Item(113, 17).Value = 20
$ws.Cells.Item(113, 18).Value = "Hortaliza"

# New row 114 - Segunda quality, newest week
$ws.Cells.Item(114, 1).Value = 11
$ws.Cells.Item(114, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(114, 3).Value = "Bíobío"
$ws.Cells.Item(114, 4).Value = 44461
$ws.Cells.Item(114, 5).Value = 8
$ws.Cells.Item(114, 6).Value = 100114013
$ws.Cells.Item(114, 7).Value = "Zanahoria"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Segunda"
$ws.Cells.Item(114, 10).Value = 300
$ws.Cells.Item(114, 11).Value = 7000
$ws.Cells.Item(114, 12).Value = 7000
$ws.Cells.Item(114, 13).Value = 7000
$ws.Cells.Item(114, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(114, 15).Value = "Región de Ñuble"
$ws.Cells.Item(114, 16).Value = 350
$ws.Cells.Item(114, 17).Value = 20
$ws.Cells.Item(114, 18).Value = "Hortaliza"
